# Update the answer values in the practice table (three-digit ÷ one-digit).
# Each data row of the table (rows 1, 5, 9, 13, 17) holds five "a÷b=c, d"
# answers; we rewrite each cell's text in place, cell by cell, so results
# land on exactly the same table positions regardless of any duplicate
# text elsewhere in the document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("581÷5=116, 1", "704÷5=140, 4", "439÷6=73, 1", "140÷8=17, 4", "383÷4=95, 3")
    5  = @("964÷3=321, 1", "731÷3=243, 2", "855÷3=285, 0", "693÷7=99, 0", "208÷3=69, 1")
    9  = @("734÷6=122, 2", "534÷5=106, 4", "948÷5=189, 3", "436÷3=145, 1", "670÷7=95, 5")
    13 = @("909÷2=454, 1", "782÷3=260, 2", "553÷9=61, 4", "991÷6=165, 1", "415÷8=51, 7")
    17 = @("434÷4=108, 2", "529÷9=58, 7", "564÷2=282, 0", "323÷6=53, 5", "880÷9=97, 7")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
